# Apply "Added NPC To The World" burndown updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the estimate/remaining-work numbers for the affected tasks.
# Row 3: Refactor to support NPC   2   -> 1
# Row 4: Implement NPC             1.5 -> 1
# Row 5: Add NPC to World          0.5 -> 0
# Row 7: NPC Movement              4   -> 3
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("C7").Value = 3

# Recalculate so the SUM(C3:C26) total in C27 (and the chart that is
# sourced from it) picks up the new figures.
$excel.CalculateFullRebuild()

# Refresh the embedded burndown chart so its cached series values match
# the recalculated totals in row 27.
$chartObj = $ws.ChartObjects(1)
$chartObj.Chart.Refresh()

# Move the active selection to reflect where work was last done.
$ws.Range("C16").Select()

$wb.Save()
